$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout change -------------------------------------------------
# A new "BUSINESS UNIT" column is inserted right after column C (becomes D),
# pushing PROD LINE / TREAD / PART# one column to the right (D->E, E->F, F->G).
# The *old* BUSINESS UNIT column (originally G, "NET PRICE" onward stays put)
# is removed so the total column count / layout (A:L) is unchanged.
$ws.Columns("G").Delete()
$ws.Columns("D").Insert()

# New header for the inserted column
$ws.Range("D1").Value = "BUSINESS UNIT"

# --- Clear the sample/demo data rows (2-6) --------------------------------
# Columns that keep their number/date/percent formatting but lose their value
$ws.Range("B2:C6").ClearContents()
$ws.Range("E2:F6").ClearContents()
$ws.Range("I2:I6").ClearContents()

# Columns that are emptied completely (value + formatting removed)
$ws.Range("A2:A6").Clear()
$ws.Range("D2:D6").Clear()
$ws.Range("G2:G6").Clear()
$ws.Range("H2:H6").Clear()
$ws.Range("J2:K6").Clear()
$ws.Range("L2:L6").Clear()

# --- Column widths (AutoFit-style re-sizing after the layout change) -----
$ws.Columns("A").ColumnWidth = 12.666666666666666
$ws.Columns("D").ColumnWidth = 14
$ws.Columns("H").ColumnWidth = 10.833333333333334
$ws.Columns("I").ColumnWidth = 13.833333333333334
$ws.Columns("J").ColumnWidth = 13.666666666666666

# --- Selection -------------------------------------------------------------
$ws.Range("G8").Select()
